$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add the new "Colorado" label next to the title (B1)
$ws.Range("B1").Value = "Colorado"

# Update the date stamp in C1 (serial date 44565 = 1/4/2022)
$ws.Range("C1").Value = 44565
